# Insert a new worksheet "2022-Q3" right after "总计" and shift the rest down.
$wb = $excel.ActiveWorkbook

# 1. Update the summary sheet ("总计") - insert a new row for 2022-Q3 at the top
#    of the data (row 2), pushing all existing quarter rows down by one.
#    (Cell-to-cell copy must use .Value2 for reads - .Value getter is not
#    reliable for this runtime when read inline; .Value2 works correctly.)
$summary = $wb.Worksheets.Item(1)
for ($r = 7; $r -ge 2; $r--) {
    $destRow = $r + 1
    $summary.Cells.Item($destRow, 1).Value = $r - 1
    $summary.Cells.Item($destRow, 2).Value = $summary.Cells.Item($r, 2).Value2
    $summary.Cells.Item($destRow, 3).Value = $summary.Cells.Item($r, 3).Value2
    $summary.Cells.Item($destRow, 4).Value = $summary.Cells.Item($r, 4).Value2
}
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 29
$summary.Range("D2").Value = 4.16

# 2. Add the new worksheet for 2022-Q3, placed right after "总计".
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Match the bold / centered / bordered look used for the header row and the
# index column (A) on the other quarter sheets in this workbook.
$headerRng = $newSheet.Range("B1:H1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

$indexRng = $newSheet.Range("A2:A30")
$indexRng.Font.Bold = $true
$indexRng.HorizontalAlignment = -4108
$indexRng.VerticalAlignment = -4160
$indexRng.Borders.LineStyle = 1

# Make B:G text columns so numeric-looking strings (fund codes, percentages)
# keep their original formatting instead of being coerced to numbers.
$newSheet.Range("B2:G30").NumberFormat = "@"

$data = @(
    @("011230","创金合信数字经济主题股票C","15.89","92.35","4.16","0.6610",6),
    @("011229","创金合信数字经济主题股票A","15.42","92.35","4.16","0.6415",6),
    @("010709","安信医药健康主题股票A","16.29","94.32","3.84","0.6255",9),
    @("010710","安信医药健康主题股票C","13.91","94.32","3.84","0.5341",9),
    @("011335","银河医药健康混合A","7.55","94.19","5.76","0.4349",4),
    @("012260","广发睿明优质企业混合A","10.34","64.35","3.53","0.3650",5),
    @("000780","鹏华医疗保健股票","6.60","81.50","3.89","0.2567",8),
    @("000339","长城医疗保健混合A","7.18","87.41","2.63","0.1888",10),
    @("519673","银河康乐股票A","1.94","93.82","4.50","0.0873",8),
    @("001060","前海开源高端装备制造灵活配置混合","0.90","80.76","7.22","0.0650",2),
    @("002662","前海开源沪港深大消费主题混合A","0.63","82.73","7.08","0.0446",4),
    @("000524","上投摩根民生需求股票","1.46","80.77","2.91","0.0425",9),
    @("002663","前海开源沪港深大消费主题混合C","0.52","82.73","7.08","0.0368",4),
    @("002515","招商丰益灵活配置混合C","1.60","39.36","2.21","0.0354",8),
    @("012261","广发睿明优质企业混合C","0.73","64.35","3.53","0.0258",5),
    @("014212","博时研究优享混合A","0.80","79.50","2.76","0.0221",10),
    @("016018","银河康乐股票C","0.35","93.82","4.50","0.0158",8),
    @("002514","招商丰益灵活配置混合A","0.60","39.36","2.21","0.0133",8),
    @("007133","嘉实长青竞争优势股票A","0.24","90.21","5.03","0.0121",9),
    @("001482","上投摩根新兴服务股票","0.38","80.81","2.99","0.0114",8),
    @("014157","国泰君安创新医药混合","0.29","79.28","3.67","0.0106",6),
    @("015201","创金合信动态平衡混合C","0.23","65.33","4.16","0.0096",2),
    @("015200","创金合信动态平衡混合A","0.16","65.33","4.16","0.0067",2),
    @("015562","长城医疗保健混合C","0.13","87.41","2.63","0.0034",10),
    @("015655","富荣医药健康混合A","0.13","82.11","1.68","0.0022",3),
    @("007134","嘉实长青竞争优势股票C","0.04","90.21","5.03","0.0020",9),
    @("014213","博时研究优享混合C","0.07","79.50","2.76","0.0019",10),
    @("015666","银河医药健康混合C","0.01","94.19","5.76","0.0006",4),
    @("015656","富荣医药健康混合C","0.01","82.11","1.68","0.0002",3)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $newSheet.Cells.Item($row, 1).Value = $i
    $newSheet.Cells.Item($row, 2).Value = $rec[0]
    $newSheet.Cells.Item($row, 3).Value = $rec[1]
    $newSheet.Cells.Item($row, 4).Value = $rec[2]
    $newSheet.Cells.Item($row, 5).Value = $rec[3]
    $newSheet.Cells.Item($row, 6).Value = $rec[4]
    $newSheet.Cells.Item($row, 7).Value = $rec[5]
    $newSheet.Cells.Item($row, 8).Value = $rec[6]
}
